# Daily auto-push update: insert a new day's record (2026/01/12, 月, 19:00, rank 24)
# right before the 2026/12/29 block, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 610 (pushes old row 610.. down to 611..)
$ws.Rows.Item(610).Insert()

# Column A holds dates as plain text (not real Excel dates) throughout the sheet,
# so force Text formatting before writing, then reset the style back to the
# sheet's default ("Normal") so no stray per-cell style is left behind.
$ws.Cells.Item(610, 1).NumberFormat = "@"
$ws.Cells.Item(610, 1).Value = "2026/01/12"
$ws.Cells.Item(610, 1).Style = "Normal"

$ws.Cells.Item(610, 2).Value = "月"
$ws.Cells.Item(610, 3).Value = 19
$ws.Cells.Item(610, 4).Value = 24
